# "Removed old postprocessing code"
# The data export used to stop at row 65 (channel 7065). With the old
# postprocessing code removed, the raw export now keeps going through
# channel 7074, and the previously-truncated last channel's count is
# corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64 (channel 7064) count corrected from 0 to 1. Text stays the same.
$ws.Range("B64").Value = 1

# Row 65 (channel 7065) is unchanged (kept for reference).

# New rows 66-74: channels 7066-7074 that postprocessing used to strip out.
$newRows = @(
    @{ Row = 66; Text = "7066"; Count = 1 },
    @{ Row = 67; Text = "7067"; Count = 2 },
    @{ Row = 68; Text = "7068"; Count = 3 },
    @{ Row = 69; Text = "7069"; Count = 2 },
    @{ Row = 70; Text = "7070"; Count = 2 },
    @{ Row = 71; Text = "7071"; Count = 2 },
    @{ Row = 72; Text = "7072"; Count = 2 },
    @{ Row = 73; Text = "7073"; Count = 1 },
    @{ Row = 74; Text = "7074"; Count = 2 }
)

foreach ($r in $newRows) {
    $cellA = $ws.Cells.Item($r.Row, 1)
    # Leading apostrophe forces text storage (numeric-looking channel id)
    # while keeping the cell on the default "Normal" style, matching the
    # rest of the sheet (s="0").
    $cellA.Value = "'" + $r.Text
    $cellA.Style = "Normal"

    $cellB = $ws.Cells.Item($r.Row, 2)
    $cellB.Value = $r.Count
}
